$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before the existing "ShiftTeamId" table column ---
# The table (Table11) currently occupies C1:F5. Inserting a plain column at
# C shifts the table's cells right, but Excel keeps the ListObject's range
# pinned to its original top-left column, so the new blank column gets
# absorbed as the table's first column. Resize the table back to D1:G5 so
# the new column becomes a normal (non-table) worksheet column, matching
# the authored change.
$lo = $ws.ListObjects.Item(1)
$ws.Columns.Item(3).Insert()
$lo.Resize($ws.Range("D1:G5"))

# Give column C roughly the same width as the neighboring ID columns.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Populate the new "KronosTimeZone" header cell (C1) ---
# Copy the header formatting (font/fill/border) from the existing
# KronosOrgJobPath header (B1) down through row 5, then set the header
# text and strip the border so it matches the un-bordered header style
# used for this new column.
$ws.Range("B1:B5").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value2 = "KronosTimeZone"
$ws.Range("C1:C5").Borders.LineStyle = -4142  # xlLineStyleNone

# --- Update the worksheet AutoFilter to span the three non-table header cells ---
$ws.AutoFilterMode = $false
$ws.Range("A1:C1").AutoFilter()

# --- Update the workbook-level _FilterDatabase defined name to match ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=" + $ws.Name + "!`$A`$1:`$C`$1"
    }
}
